# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Each coin's row keeps its rank (column A) and link; only the fields the
# scraper re-read -- Coin/Link/Price/Volume(1h) -- are rewritten in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price text that parses as a plain number (e.g. "227.02") must be entered
# with a leading apostrophe so Excel stores it as text, matching the source
# data's inline-string cells, instead of silently converting it to a Double.

# Row 2
$ws.Cells.Item(2, 4).Value = '34.115.70'
$ws.Cells.Item(2, 5).Value = '  +0.09%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.789.77'
$ws.Cells.Item(3, 5).Value = '  -0.16%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.07%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''227.02'
$ws.Cells.Item(5, 5).Value = '  +1.08%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''0.545'
$ws.Cells.Item(6, 5).Value = '  -0.79%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.05%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''32.33'
$ws.Cells.Item(8, 5).Value = '  -0.22%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +4.09%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.0688'
$ws.Cells.Item(10, 5).Value = '  -2.67%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.28%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '2.049.33'
$ws.Cells.Item(12, 5).Value = '  -0.09%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''11.47'
$ws.Cells.Item(13, 5).Value = '  +6.10%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.800.95'
$ws.Cells.Item(14, 5).Value = '  +0.25%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''0.623'

# Row 16
$ws.Cells.Item(16, 4).Value = '34.106.28'
$ws.Cells.Item(16, 5).Value = '  +0.14%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.50%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''67.98'
$ws.Cells.Item(18, 5).Value = '  +0.04%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''244.12'
$ws.Cells.Item(19, 5).Value = '  +0.30%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0777'
$ws.Cells.Item(20, 5).Value = '  -0.91%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.07%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''10.92'
$ws.Cells.Item(22, 5).Value = '  +2.39%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''4.10'
$ws.Cells.Item(23, 5).Value = '  +0.58%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -2.05%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''161.56'
$ws.Cells.Item(25, 5).Value = '  +1.70%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +2.41%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''16.26'
$ws.Cells.Item(27, 5).Value = '  +0.20%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +1.26%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.17%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.89%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.21%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.10%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''3.62'
$ws.Cells.Item(33, 5).Value = '  +3.79%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +1.41%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '1.404.97'
$ws.Cells.Item(35, 5).Value = '  +1.40%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''0.654'
$ws.Cells.Item(36, 5).Value = '  +1.24%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''1.04'
$ws.Cells.Item(37, 5).Value = '  -0.36%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +2.32%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''2.34'
$ws.Cells.Item(39, 5).Value = '  +8.01%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''80.08'
$ws.Cells.Item(40, 5).Value = '  +1.25%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.50%  '

# Row 42 - ARBITRUM and MXToken swapped rank order with this refresh, so
# Coin/Link/Price/Volume are rewritten here (row 42 now holds MXToken)
$ws.Cells.Item(42, 2).Value = 'MXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(42, 4).Value = '''2.72'
$ws.Cells.Item(42, 5).Value = '  +0.82%  '

# Row 43 - ... and row 43 now holds ARBITRUM
$ws.Cells.Item(43, 2).Value = 'ARBITRUM'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(43, 4).Value = '''0.922'
$ws.Cells.Item(43, 5).Value = '  +0.84%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''13.41'
$ws.Cells.Item(44, 5).Value = '  +12.33%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +1.25%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''6.07'
$ws.Cells.Item(46, 5).Value = '  +3.94%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.0508'
$ws.Cells.Item(47, 5).Value = '  +2.72%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +2.67%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''107.21'
$ws.Cells.Item(49, 5).Value = '  +0.16%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '1.950.77'
$ws.Cells.Item(50, 5).Value = '  +0.06%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.09%  '
